# ch3 processes-drivers
#
# Adds a new worksheet "历史时期特征划分" (historical-period feature
# breakdown) summarising time-span / period-division / main-characteristics
# for the Yellow River basin, and updates the active selections on the
# pre-existing sheets accordingly.

$wb = $excel.ActiveWorkbook

# --- 1. Move the selection on sheet 1 (existing data-source sheet) -------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("K14").Select()

# --- 2. Add the new worksheet as the LAST tab -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "历史时期特征划分"

# Column widths (characters), matching the other summary tables in the book
$ws3.Columns.Item(1).ColumnWidth = 17
$ws3.Columns.Item(2).ColumnWidth = 16.33
$ws3.Columns.Item(3).ColumnWidth = 35.83

# --- 3. Header row ---------------------------------------------------------
$ws3.Range("A2").Value = "时间跨度"
$ws3.Range("B2").Value = "时段划分"
$ws3.Range("C2").Value = "主要特征"

# --- 4. Data rows ------------------------------------------------------------
$ws3.Range("A3").Value = "200BC-400AD"
$ws3.Range("B3").Value = "数据不可信时段"
$ws3.Range("C3").Value = "各数据集在此时期的可信度均偏低"

$ws3.Range("A4").Value = "400-900AD"
$ws3.Range("B4").Value = "CDP1前期"
$ws3.Range("C4").Value = "没有明显的驱动因素"

$ws3.Range("A5").Value = "900-1100AD"
$ws3.Range("B5").Value = "CDP1时期"
$ws3.Range("C5").Value = "气候驱动与低水平的人类活动驱动时期"

$ws3.Range("A6").Value = "1100-1350AD"
$ws3.Range("B6").Value = "CDP1后期"
$ws3.Range("C6").Value = "没有明显的驱动因素"

$ws3.Range("A7").Value = "1350-1700AD"
$ws3.Range("B7").Value = "CDP2前期"
$ws3.Range("C7").Value = "人类活动驱动时期"

$ws3.Range("A8").Value = "1700-1900AD"
$ws3.Range("B8").Value = "CDP2时期"
$ws3.Range("C8").Value = "气候驱动与人类活动共同驱动时期"

$ws3.Range("A9").Value = "1900-2000AD"
$ws3.Range("B9").Value = "HDP2时期"
$ws3.Range("C9").Value = "人口迅速增长的人类活动强烈驱动期"

# --- 5. Formatting, reusing the exact look used by the other tables in ---
#        this workbook: header row -> top+bottom rule, last row -> bottom
#        rule. Copy the formatting from the equivalent rows on sheet 1
#        so the same style entries are (re)used instead of minting new
#        near-duplicate ones.
$ws1.Range("A2").Copy()
$ws3.Range("A2:C2").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("A9").Copy()
$ws3.Range("A9:C9").PasteSpecial(-4122)  # xlPasteFormats

$ws3.Range("A2:C9").VerticalAlignment = -4108  # xlCenter
$excel.CutCopyMode = 0

# --- 6. Selection / active cell on the new sheet, matches the diff --------
$ws3.Range("D9").Select()
